$d = $word.ActiveDocument

# The table to update is the last table in the document (the "Sprint 4" log
# table whose second row is still blank). Locate it, remember where it
# starts, delete it, then re-insert the fully updated OOXML for the table
# (filled-in row + recalculated tblGrid column widths) at that position.
$t = $d.Tables.Item($d.Tables.Count)
$tblStart = $t.Range.Start
$t.Delete()

$newTableXml = '<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblpPr w:leftFromText="141" w:rightFromText="141" w:vertAnchor="text" w:horzAnchor="margin" w:tblpXSpec="center" w:tblpY="277"/><w:tblW w:w="10484" w:type="dxa"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="1394"/><w:gridCol w:w="3432"/><w:gridCol w:w="1696"/><w:gridCol w:w="1646"/><w:gridCol w:w="1083"/><w:gridCol w:w="1233"/></w:tblGrid><w:tr w:rsidR="0078770F" w14:paraId="796FAC06" w14:textId="77777777" w:rsidTr="00A2076E"><w:tc><w:tcPr><w:tcW w:w="1397" w:type="dxa"/><w:shd w:val="clear" w:color="auto" w:fill="D9D9D9" w:themeFill="background1" w:themeFillShade="D9"/></w:tcPr><w:p w14:paraId="2ACC21AC" w14:textId="77777777" w:rsidR="0078770F" w:rsidRDefault="0078770F" w:rsidP="00A2076E"><w:pPr><w:rPr><w:lang w:val="x-none" w:eastAsia="x-none"/></w:rPr></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3491" w:type="dxa"/><w:shd w:val="clear" w:color="auto" w:fill="D9D9D9" w:themeFill="background1" w:themeFillShade="D9"/></w:tcPr><w:p w14:paraId="74349AC6" w14:textId="77777777" w:rsidR="0078770F" w:rsidRDefault="0078770F" w:rsidP="00A2076E"><w:pPr><w:rPr><w:lang w:val="x-none" w:eastAsia="x-none"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="x-none" w:eastAsia="x-none"/></w:rPr><w:t>Tarea</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1616" w:type="dxa"/><w:shd w:val="clear" w:color="auto" w:fill="D9D9D9" w:themeFill="background1" w:themeFillShade="D9"/></w:tcPr><w:p w14:paraId="21228F39" w14:textId="77777777" w:rsidR="0078770F" w:rsidRDefault="0078770F" w:rsidP="00A2076E"><w:pPr><w:rPr><w:lang w:val="x-none" w:eastAsia="x-none"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="x-none" w:eastAsia="x-none"/></w:rPr><w:t>Etiqueta</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1655" w:type="dxa"/><w:shd w:val="clear" w:color="auto" w:fill="D9D9D9" w:themeFill="background1" w:themeFillShade="D9"/></w:tcPr><w:p w14:paraId="45DD49C2" w14:textId="77777777" w:rsidR="0078770F" w:rsidRDefault="0078770F" w:rsidP="00A2076E"><w:pPr><w:rPr><w:lang w:val="x-none" w:eastAsia="x-none"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="x-none" w:eastAsia="x-none"/></w:rPr><w:t>Estado</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1083" w:type="dxa"/><w:shd w:val="clear" w:color="auto" w:fill="D9D9D9" w:themeFill="background1" w:themeFillShade="D9"/></w:tcPr><w:p w14:paraId="55637F54" w14:textId="77777777" w:rsidR="0078770F" w:rsidRDefault="0078770F" w:rsidP="00A2076E"><w:pPr><w:rPr><w:lang w:val="x-none" w:eastAsia="x-none"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="x-none" w:eastAsia="x-none"/></w:rPr><w:t>Tiempo estimado</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1242" w:type="dxa"/><w:shd w:val="clear" w:color="auto" w:fill="D9D9D9" w:themeFill="background1" w:themeFillShade="D9"/></w:tcPr><w:p w14:paraId="651889E9" w14:textId="77777777" w:rsidR="0078770F" w:rsidRDefault="0078770F" w:rsidP="00A2076E"><w:pPr><w:rPr><w:lang w:val="x-none" w:eastAsia="x-none"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="x-none" w:eastAsia="x-none"/></w:rPr><w:t>Tiempo real</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="0078770F" w14:paraId="49246109" w14:textId="77777777" w:rsidTr="00A2076E"><w:tc><w:tcPr><w:tcW w:w="1397" w:type="dxa"/><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="4C275019" w14:textId="20BF9BBF" w:rsidR="0078770F" w:rsidRDefault="0078770F" w:rsidP="00A2076E"><w:pPr><w:rPr><w:lang w:val="x-none" w:eastAsia="x-none"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="x-none" w:eastAsia="x-none"/></w:rPr><w:t>16/04/2024</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3491" w:type="dxa"/></w:tcPr><w:p w14:paraId="1B38BAE8" w14:textId="71C6BBC9" w:rsidR="0078770F" w:rsidRDefault="0078770F" w:rsidP="00A2076E"><w:pPr><w:rPr><w:lang w:val="x-none" w:eastAsia="x-none"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="x-none" w:eastAsia="x-none"/></w:rPr><w:t>Documentación del sprint 3 y cambios en las tablas y gráficos</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1616" w:type="dxa"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/></w:tcPr><w:p w14:paraId="6A67185E" w14:textId="701F0035" w:rsidR="0078770F" w:rsidRPr="004E3232" w:rsidRDefault="0078770F" w:rsidP="00A2076E"><w:pPr><w:jc w:val="center"/><w:rPr><w:lang w:val="x-none" w:eastAsia="x-none"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="x-none" w:eastAsia="x-none"/></w:rPr><w:t>Documentation</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1655" w:type="dxa"/></w:tcPr><w:p w14:paraId="101EAFEF" w14:textId="7057EA8F" w:rsidR="0078770F" w:rsidRDefault="0078770F" w:rsidP="00A2076E"><w:pPr><w:jc w:val="center"/><w:rPr><w:lang w:val="x-none" w:eastAsia="x-none"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="x-none" w:eastAsia="x-none"/></w:rPr><w:t>Completado</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1083" w:type="dxa"/></w:tcPr><w:p w14:paraId="045018E1" w14:textId="77CE7325" w:rsidR="0078770F" w:rsidRDefault="0078770F" w:rsidP="00A2076E"><w:pPr><w:jc w:val="center"/><w:rPr><w:lang w:val="x-none" w:eastAsia="x-none"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="x-none" w:eastAsia="x-none"/></w:rPr><w:t>0.5</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1242" w:type="dxa"/></w:tcPr><w:p w14:paraId="0FD713EF" w14:textId="2B2EB63C" w:rsidR="0078770F" w:rsidRDefault="0078770F" w:rsidP="00A2076E"><w:pPr><w:jc w:val="center"/><w:rPr><w:lang w:val="x-none" w:eastAsia="x-none"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="x-none" w:eastAsia="x-none"/></w:rPr><w:t>0.5</w:t></w:r></w:p></w:tc></w:tr></w:tbl>'

$insertRange = $d.Range($tblStart, $tblStart)
$insertRange.InsertXML($newTableXml)
